# Renamed few transcripts. Updated the DataSheet.
# Column D ("Speaker") values "RT1" get shortened to "T" and the single
# "A student" entry becomes "S" for the rows listed below.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$rt1Rows = @(2,3,4,6,8,9,10,12,13,14,17,18,19,25,30,31,32,36,37,38,50,52,54,55,59,60,61,62,63,64,65,71,73,79,81,82,83,85,87,89,90,91,96,98,99,100,104,106,117,119,120,128,131,132,133,134,135,136,138,139,140,141,142,143,144,145,146,147,148,150,153,158,159,160,162,163,164,165,166,167,169,170,172,174,181,182,183,184,191,193,194,195,196,197,198,199,200,201,202,203,206,207,211,212,213,214,215,216,217,219,222,224,225,233,234,235,241,243,244,246,248)

foreach ($r in $rt1Rows) {
    $ws.Cells.Item($r, 4).Value = "T"
}

# Row 22: "A student" -> "S"
$ws.Cells.Item(22, 4).Value = "S"
